$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 3335.322
$ws.Range("I18").Value = 3444.9812
$ws.Range("J18").Value = 2366.6667
$ws.Range("K18").Value = 3444.9812
$ws.Range("L18").Value = 2366.6667
$ws.Range("M18").Value = -3160.9812
$ws.Range("N18").Value = -2934.6667

$ws.Range("H28").Value = 3870.6667
$ws.Range("I28").Value = 4334.5
$ws.Range("K28").Value = 4334.5
$ws.Range("M28").Value = -3849.5

$ws.Range("H33").Value = 677.3889
$ws.Range("I33").Value = 249.78572
$ws.Range("J33").Value = 2174
$ws.Range("K33").Value = 249.78572
$ws.Range("L33").Value = 2174
$ws.Range("M33").Value = -20.78572
$ws.Range("N33").Value = -2632

$ws.Range("H46").Value = 84306.92
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 84306.92
$ws.Range("K46").Value = 0
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -253158.76

$ws.Range("H60").Value = 84306.92
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 84306.92
$ws.Range("K60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -253888.76

$ws.Range("H80").Value = 158517.92
$ws.Range("I80").Value = 286541.16
$ws.Range("J80").Value = 9157.5
$ws.Range("K80").Value = 859623.48
$ws.Range("L80").Value = 27472.5
$ws.Range("M80").Value = -858625.48
$ws.Range("N80").Value = -29468.5

$ws.Range("H83").Value = 158517.92
$ws.Range("I83").Value = 286541.16
$ws.Range("J83").Value = 9157.5
$ws.Range("K83").Value = 2578870.44
$ws.Range("L83").Value = 82417.5
$ws.Range("M83").Value = -2573878.44
$ws.Range("N83").Value = -92401.5

$ws.Range("H116").Value = 444991.56
$ws.Range("I116").Value = 745801.9
$ws.Range("J116").Value = 5345.6924
$ws.Range("K116").Value = 745801.9
$ws.Range("L116").Value = 5345.6924
$ws.Range("M116").Value = -742359.9
$ws.Range("N116").Value = -12229.6924

$ws.Range("H132").Value = 4891.8696
$ws.Range("I132").Value = 4771.7617
$ws.Range("K132").Value = 14315.2851
$ws.Range("M132").Value = -11785.2851

$ws.Range("H138").Value = 4762
$ws.Range("J138").Value = 5538.237
$ws.Range("L138").Value = 16614.711
$ws.Range("N138").Value = -26894.711

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1551.125
$ws.Range("I63").Value = 1551.125
$ws.Range("K63").Value = 1551.125
$ws.Range("M63").Value = -865.125

$ws.Range("H66").Value = 1551.125
$ws.Range("I66").Value = 1551.125
$ws.Range("K66").Value = 7755.625
$ws.Range("M66").Value = -4323.625

$ws.Range("H97").Value = 8515.25
$ws.Range("I97").Value = 9727.615
$ws.Range("K97").Value = 9727.615
$ws.Range("M97").Value = -9231.615

$ws.Range("H122").Value = 2204528.5
$ws.Range("I122").Value = 3623.1428
$ws.Range("K122").Value = 10869.4284
$ws.Range("M122").Value = -8419.428400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2335.1738
$ws.Range("I94").Value = 1937.1052
$ws.Range("K94").Value = 1937.1052
$ws.Range("M94").Value = -1486.1052

$ws.Range("H105").Value = 2570.6155
$ws.Range("I105").Value = 2098.625
$ws.Range("K105").Value = 2098.625
$ws.Range("M105").Value = -351.625

$ws.Range("H134").Value = 5268.6665
$ws.Range("I134").Value = 5024.8887
$ws.Range("K134").Value = 15074.6661
$ws.Range("M134").Value = -12539.6661

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 709
$ws.Range("I22").Value = 736.6875
$ws.Range("J22").Value = 677.3570999999999
$ws.Range("K22").Value = 736.6875
$ws.Range("L22").Value = 677.3570999999999
$ws.Range("M22").Value = -386.6875
$ws.Range("N22").Value = -1377.3571

$ws.Range("H58").Value = 2101.4443
$ws.Range("I58").Value = 1479.08
$ws.Range("K58").Value = 1479.08
$ws.Range("M58").Value = -1276.08

$ws.Range("H122").Value = 758.8461
$ws.Range("I122").Value = 758.8461
$ws.Range("K122").Value = 2276.5383
$ws.Range("M122").Value = 173.4616999999998

$ws.Range("H132").Value = 17285.03
$ws.Range("I132").Value = 6028.706
$ws.Range("J132").Value = 400000
$ws.Range("K132").Value = 18086.118
$ws.Range("L132").Value = 1200000
$ws.Range("M132").Value = -15556.118
$ws.Range("N132").Value = -1205060

$ws.Range("H134").Value = 3685537.2
$ws.Range("I134").Value = 4176009
$ws.Range("K134").Value = 12528027
$ws.Range("M134").Value = -12525492

$ws.Range("H136").Value = 2101.4443
$ws.Range("I136").Value = 1479.08
$ws.Range("K136").Value = 4437.24
$ws.Range("M136").Value = -1887.24

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4568490
$ws.Range("I4").Value = 2286404
$ws.Range("K4").Value = 6859212
$ws.Range("M4").Value = -6859100

$ws.Range("H5").Value = 346549.38
$ws.Range("J5").Value = 770842.4
$ws.Range("L5").Value = 2312527.2
$ws.Range("N5").Value = -2312751.2

$ws.Range("H37").Value = 61546.89
$ws.Range("J37").Value = 61546.89
$ws.Range("L37").Value = 184640.67
$ws.Range("N37").Value = -184864.67

$ws.Range("H39").Value = 4928.5713
$ws.Range("I39").Value = 4916.6665
$ws.Range("J39").Value = 5000
$ws.Range("K39").Value = 14749.9995
$ws.Range("L39").Value = 15000
$ws.Range("M39").Value = -14455.9995
$ws.Range("N39").Value = -15588

$ws.Range("H46").Value = 2118.6155
$ws.Range("I46").Value = 498.66666
$ws.Range("J46").Value = 3507.1428
$ws.Range("K46").Value = 1495.99998
$ws.Range("L46").Value = 10521.4284
$ws.Range("M46").Value = -1404.99998
$ws.Range("N46").Value = -10703.4284

$ws.Range("H51").Value = 668359.8
$ws.Range("I51").Value = 909995.4399999999
$ws.Range("J51").Value = 3861.75
$ws.Range("K51").Value = 2729986.32
$ws.Range("L51").Value = 11585.25
$ws.Range("M51").Value = -2729526.32
$ws.Range("N51").Value = -12505.25

$ws.Range("H58").Value = 1400.3334
$ws.Range("I58").Value = 1400.3334
$ws.Range("K58").Value = 4201.0002
$ws.Range("M58").Value = -4073.0002

$ws.Range("H122").Value = 1366.5264
$ws.Range("J122").Value = 1652.6364
$ws.Range("L122").Value = 14873.7276
$ws.Range("N122").Value = -19773.7276

$ws.Range("H131").Value = 47620252
$ws.Range("I131").Value = 100000584
$ws.Range("J131").Value = 1766.3636
$ws.Range("K131").Value = 300001752
$ws.Range("L131").Value = 5299.0908
$ws.Range("M131").Value = -299996712
$ws.Range("N131").Value = -15379.0908

$ws.Range("H135").Value = 346549.38
$ws.Range("J135").Value = 770842.4
$ws.Range("L135").Value = 6937581.600000001
$ws.Range("N135").Value = -6942651.600000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2459.8809
$ws.Range("I70").Value = 2320.1
$ws.Range("K70").Value = 2320.1
$ws.Range("M70").Value = -2050.1

$ws.Range("H73").Value = 2459.8809
$ws.Range("I73").Value = 2320.1
$ws.Range("K73").Value = 2320.1
$ws.Range("M73").Value = -1384.1

$ws.Range("H97").Value = 7765.5
$ws.Range("I97").Value = 8882.120000000001
$ws.Range("J97").Value = 3777.5715
$ws.Range("K97").Value = 8882.120000000001
$ws.Range("L97").Value = 3777.5715
$ws.Range("M97").Value = -8386.120000000001
$ws.Range("N97").Value = -4769.5715

$ws.Range("H107").Value = 586.41174
$ws.Range("J107").Value = 870.25
$ws.Range("L107").Value = 870.25
$ws.Range("N107").Value = -4710.25

$ws.Range("H122").Value = 29627.334
$ws.Range("I122").Value = 33682.57
$ws.Range("K122").Value = 101047.71
$ws.Range("M122").Value = -98597.70999999999

$ws.Range("H132").Value = 3759.8948
$ws.Range("I132").Value = 2624.25
$ws.Range("K132").Value = 7872.75
$ws.Range("M132").Value = -5342.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 31643.234
$ws.Range("I7").Value = 38597.152
$ws.Range("K7").Value = 38597.152
$ws.Range("M7").Value = -38485.152

$ws.Range("H16").Value = 3427
$ws.Range("I16").Value = 2977.5
$ws.Range("K16").Value = 2977.5
$ws.Range("M16").Value = -2807.5

$ws.Range("H126").Value = 31643.234
$ws.Range("I126").Value = 38597.152
$ws.Range("K126").Value = 115791.456
$ws.Range("M126").Value = -113321.456

$ws.Range("H132").Value = 1450046.9
$ws.Range("I132").Value = 2393452.2
$ws.Range("J132").Value = 3492.2666
$ws.Range("K132").Value = 7180356.600000001
$ws.Range("L132").Value = 10476.7998
$ws.Range("M132").Value = -7177826.600000001
$ws.Range("N132").Value = -15536.7998

$ws.Range("H136").Value = 12711.429
$ws.Range("I136").Value = 1490
$ws.Range("K136").Value = 4470
$ws.Range("M136").Value = -1920

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 28771.953
$ws.Range("I122").Value = 2622.8333
$ws.Range("K122").Value = 7868.499899999999
$ws.Range("M122").Value = -5418.499899999999

$ws.Range("H132").Value = 9499.270500000001
$ws.Range("I132").Value = 10386.65
$ws.Range("K132").Value = 31159.95
$ws.Range("M132").Value = -28629.95

$ws.Range("H136").Value = 1850.2122
$ws.Range("I136").Value = 1253.6129
$ws.Range("J136").Value = 11097.5
$ws.Range("K136").Value = 3760.8387
$ws.Range("L136").Value = 33292.5
$ws.Range("M136").Value = -1210.8387
$ws.Range("N136").Value = -38392.5
